# Fix 2-column vs 3-column confusion in the Illinois Birth Weights workbook.
# Original layout: A = "Born in US", B = "Born in Africa", C = Comments block.
# New layout:       A = "Black_US", B = "Black_Africa", C = "White_US" (new
#                    data column), D = Comments block (shifted right by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Insert a new column at C, shifting the old Comments column (C) to D.
$ws.Range("C:C").Insert()

# 2) Rename the two existing headers and add the new column's header.
$ws.Range("A1").Value = "Black_US"
$ws.Range("B1").Value = "Black_Africa"
$ws.Range("C1").Value = "White_US"

# 3) Fill in the new "White_US" data column (rows 2-45).
$whiteUs = @(4429,3191,3712,3399,2638,3946,3173,2926,2303,3885,3208,2969,2948,2270,3172,2318,2456,3661,3854,3122,3666,4414,3490,3871,2679,2850,2852,3316,3596,2719,4448,3043,2709,3695,3583,2867,4056,3342,3124,4281,3839,3458,3931,4322)

for ($i = 0; $i -lt $whiteUs.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $whiteUs[$i]
}

# 4) Update the Comments text in column D to describe the three categories.
$ws.Range("D10").Value = "of children born in Illinois to mothers who fall into"
$ws.Range("D11").Value = "one of the following categories:"
$ws.Range("D12").Value = "(1) Black, born in the United States (Black_US)"
$ws.Range("D13").Value = "(2) Black, born in Africa (Black_Africa), or"
$ws.Range("D14").Value = "(3) White, born in the United States (White_US)."

# 5) Column D previously stopped at row 13 (no cells below). The new file
#    carries the Comments column down to row 45 as explicit empty strings,
#    matching the extended White_US column.
for ($row = 15; $row -le 45; $row++) {
    $ws.Cells.Item($row, 4).Value = ""
}
